$wb = $excel.ActiveWorkbook

# ---- Sheet "Add Panels" (first sheet) ----
$ws1 = $wb.Worksheets.Item("Add Panels")

# Insert a new column before column J ("Minimum Battery Value") for the new
# "Alarm Hours" column, cloning column I's formatting (shifts old J->K, old K->L).
$ws1.Columns("I").Copy()
$ws1.Columns("J").Insert()

# Delete the old "FIRECLASS 64-2" data row (row 8); the FC702S row (row 9)
# shifts up to become row 8.
$ws1.Rows("8").Delete()

# Header / description updates
$ws1.Range("B3").Value = "verifyStandyByAlarmHourAndBatteryFactor"
$ws1.Range("B4").Value = "NGC-1928/T962 OR TC-71697"
$ws1.Range("J7").Value = "Alarm Hours"

# New Alarm Hours value for the remaining data row (restore the quoted-text
# style the cell inherited from column I after the Value write resets it).
$ws1.Range("J8").Value = 0.75
$ws1.Range("I8").Copy()
$ws1.Range("J8").PasteSpecial(-4122)

# Updated Minimum Battery Value for the remaining data row
$ws1.Range("K8").Value = 37.799999999999997

$ws1.Range("I3").Select()

# ---- Sheet "Sheet1" (second sheet) ----
$ws2 = $wb.Worksheets.Item("Sheet1")

# Clone row 8 (keeps all formatting/styles) into row 9, then edit the values
# that differ for the new "FC64-2" panel entry.
$ws2.Range("A8:K8").Copy($ws2.Range("A9:K9"))

$ws2.Range("A9").Value = "FC64-2"
$ws2.Range("F9").Value = 20.39
$ws2.Range("J9").Value = 25.81

$ws2.Range("A9:XFD9").Select()

$ws1.Activate()
